$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) numeric-looking text updates (use leading apostrophe to force text) ---
$ws.Range("D2").Value  = "'289.70"
$ws.Range("D3").Value  = "'21.31"
$ws.Range("D5").Value  = "'0.06389"
$ws.Range("D6").Value  = "'3.601"
$ws.Range("D7").Value  = "'1.581"
$ws.Range("D8").Value  = "'6.587"
$ws.Range("D9").Value  = "'0.8273"
$ws.Range("D10").Value = "'0.01428"
$ws.Range("D11").Value = "'0.1686"
$ws.Range("D12").Value = "'0.08803"
$ws.Range("D13").Value = "'0.03683"
$ws.Range("D14").Value = "'0.03208"
$ws.Range("D15").Value = "'0.09192"
$ws.Range("D16").Value = "'3.703"
$ws.Range("D17").Value = "'0.001663"
$ws.Range("D18").Value = "'0.04757"
$ws.Range("D19").Value = "'0.006158"
$ws.Range("D20").Value = "'0.006304"
$ws.Range("D23").Value = "'3.780"

# --- Row 27 text correction ---
$ws.Range("E27").Value = "26AAXTokenAABBestin24h"

$ws.Range("D28").Value = "'0.0002708"
$ws.Range("D40").Value = "'0.04837"
$ws.Range("D41").Value = "'0.007160"

# --- Rows 42 and 43 swap (CEJI <-> BKEXToken) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1121"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003465"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.01181"
$ws.Range("D45").Value = "'0.00006912"
$ws.Range("D47").Value = "'0.9339"

$ws.Range("D48").Value = "'0.008289"
$ws.Range("E48").Value = "47BOLOBOLO"

$ws.Range("D49").Value = "'0.00001901"
$ws.Range("D50").Value = "'0.01241"
